$d = $word.ActiveDocument

# Locate the paragraph that holds "Srbija do Tokija" - the three new
# paragraphs (two blank, one with "Ideeeemmmmoooooo") are inserted right
# after it, before the blank paragraphs that already followed it.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Srbija do Tokija*") {
        $anchor = $p
        break
    }
}

$anchorIndex = $anchor.Range.Start

# Insert paragraph #1 (blank) right after the anchor paragraph.
$anchor.Range.InsertParagraphAfter()

# Insert paragraph #2 (blank) right after the one we just created.
$p2 = $d.Paragraphs.Item($anchor.Index + 1)
$p2.Range.InsertParagraphAfter()

# Insert paragraph #3 (blank, will hold the new text) right after that.
$p3 = $d.Paragraphs.Item($anchor.Index + 2)
$p3.Range.InsertParagraphAfter()

# Fill the third new paragraph with the new sentence.
$p4 = $d.Paragraphs.Item($anchor.Index + 3)
$p4.Range.Text = "Ideeeemmmmoooooo"

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
